$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATAPNC")

# New contact numbers (column B) and vehicle numbers (column E) for rows 2-13
$contactNumbers = @(123456790144,123456790145,123456790146,123456790147,123456790148,123456790149,123456790150,123456790151,123456790152,123456790153,123456790154,123456790155)
$vehicleNumbers = @("TAA255","TAA256","TAA257","TAA258","TAA259","TAA260","TAA261","TAA262","TAA263","TAA264","TAA265","TAA266")

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $contactNumbers[$i]
    $ws.Cells.Item($row, 5).Value = $vehicleNumbers[$i]
}

# Rows 10-13 (providerList, column Y) change from "TOWING XYZ" to "ATA Baterikau"
for ($row = 10; $row -le 13; $row++) {
    $ws.Range("Y" + $row).Value = "ATA Baterikau"
}
